$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new value would otherwise be
# auto-parsed as a number by Excel, so they stay text like the source data.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.083.86"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "1.651.47"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").Value = "218.22"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").Value = "0.5292"
$ws.Range("E6").Value = "  +1.37%  "
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").Value = "0.2611"
$ws.Range("E8").Value = "  -2.06%  "
$ws.Range("D9").Value = "0.06312"
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").Value = "20.37"
$ws.Range("E10").Value = "  -3.42%  "
$ws.Range("D11").Value = "0.07740"
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("D12").Value = "4.473"
$ws.Range("E12").Value = "  +0.95%  "
$ws.Range("D13").Value = "1.653.79"
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("D14").Value = "0.5456"
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("D15").Value = "0.0₅8107"
$ws.Range("E15").Value = "  -1.58%  "
$ws.Range("D16").Value = "65.09"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").Value = "26.099.44"
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").Value = "4.550"
$ws.Range("E19").Value = "  -2.50%  "
$ws.Range("D20").Value = "193.61"
$ws.Range("E20").Value = "  +0.25%  "
$ws.Range("D21").Value = "10.04"
$ws.Range("E21").Value = "  -1.07%  "
$ws.Range("D22").Value = "5.989"
$ws.Range("E22").Value = "  -1.87%  "
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("D24").Value = "140.24"
$ws.Range("E24").Value = "  +1.17%  "
$ws.Range("D25").Value = "0.1239"
$ws.Range("E25").Value = "  -0.29%  "
$ws.Range("D26").Value = "7.238"
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("E28").Value = "  +1.01%  "
$ws.Range("D29").Value = "0.05909"
$ws.Range("E29").Value = "  -1.81%  "
$ws.Range("E30").Value = "  -0.46%  "
$ws.Range("D31").Value = "3.499"
$ws.Range("E31").Value = "  -2.15%  "
$ws.Range("D32").Value = "3.238"
$ws.Range("E32").Value = "  -2.81%  "
$ws.Range("D33").Value = "1.545"
$ws.Range("E33").Value = "  -6.39%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "0.9427"
$ws.Range("E35").Value = "  -4.03%  "
$ws.Range("D36").Value = "2.756"
$ws.Range("E36").Value = "  -0.85%  "
$ws.Range("D37").Value = "0.5653"
$ws.Range("E37").Value = "  -4.62%  "
$ws.Range("D38").Value = "0.01606"
$ws.Range("E38").Value = "  +1.23%  "
$ws.Range("D39").Value = "5.841"
$ws.Range("E39").Value = "  -1.91%  "
$ws.Range("D40").Value = "0.8439"
$ws.Range("E40").Value = "  -2.24%  "
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "100.71"
$ws.Range("E42").Value = "  +1.05%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.006.69"
$ws.Range("E43").Value = "  -3.27%  "
$ws.Range("D44").Value = "1.797.63"
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("D45").Value = "56.87"
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("E46").Value = "  -5.01%  "
$ws.Range("E47").Value = "  +0.33%  "
$ws.Range("E48").Value = "  +1.55%  "
$ws.Range("D49").Value = "1.480"
$ws.Range("E49").Value = "  +1.43%  "
$ws.Range("E50").Value = "  -0.60%  "
$ws.Range("D51").Value = "7.796"
$ws.Range("E51").Value = "  -3.59%  "
